# The document has two headers (BTec_Logo-Orange, id 1 & 3) and two
# footers (PearsonLogo.png, id 2 & 4) each carrying a single inline
# picture. This rename swaps the "friendly" picture names:
#   BTec logo pictures:    image1.jpg -> image2.jpg
#   Pearson logo pictures: image2.png -> image1.png
#
# Renaming is routed through Selection (.Range.Select() then
# $word.Selection.InlineShapes) rather than the HeaderFooter.Range
# InlineShapes collection directly, since setting .Name on that
# collection's items doesn't reliably commit for footer stories.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-LogoPicture($story, $newName) {
    for ($j = 1; $j -le $story.Range.InlineShapes.Count; $j++) {
        $shp = $story.Range.InlineShapes.Item($j)
        $shp.Range.Select()
        $word.Selection.InlineShapes.Item(1).Name = $newName
    }
}

# Headers - BTec_Logo-Orange pictures: image1.jpg -> image2.jpg
for ($i = 1; $i -le $sec.Headers.Count; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        Rename-LogoPicture $hdr "image2.jpg"
    }
}

# Footers - Pearson logo pictures: image2.png -> image1.png
for ($i = 1; $i -le $sec.Footers.Count; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        Rename-LogoPicture $ftr "image1.png"
    }
}

Write-Output "Renamed header/footer logo pictures"
